# Timing.xlsx - 2023 Advent of Code update
# Adds the Day 20 timing result: Part I took 94 ms.
# (Part II for day 20 was never solved, so column E stays blank,
#  matching the rest of the unfinished days below it.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate() | Out-Null

# Day 20 (row 25) Part I timing, in ms. G25 (=D25+E25), the Average/Total
# rows (D32/G32/D34/G34) and the bar chart all recalculate automatically.
$ws.Range("D25").Value = 94

# Leave the selection/cursor where the author last left it.
$ws.Range("F24").Select() | Out-Null
